$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colValues = @{
    "D"  = 0.133
    "E"  = 0.29
    "K"  = 14.9
    "L"  = 0.2564543889845095
    "M"  = 0.047
    "N"  = 0.0002379746835443038
    "O"  = 0.003154362416107382
    "P"  = 0.047
    "Q"  = 0.0002379746835443038
    "R"  = 0.003154362416107382
    "U"  = 367.9
    "V"  = 1.862784810126582
    "W"  = 0.07457457457457457
    "X"  = 0.06088784682880495
    "Y"  = 0.01368672774576962
    "Z"  = -0.423469387755102
    "AA" = -0
    "AB" = 0.06009687171743465
    "AC" = -0.06009687171743465
    "AD" = 12.1
    "AF" = 12.1
    "AG" = -355.8
    "AH" = 0.05772900763358779
    "AI" = 0.05339805825242718
    "AJ" = 2.247631080227417
    "AK" = 2.518046709129512
}

foreach ($row in 2..3) {
    foreach ($col in $colValues.Keys) {
        $ws.Range("$col$row").Value = $colValues[$col]
    }
}
